$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (for 2021-10-14) was inserted between the
# existing row 67 and what used to be row 68, pushing every subsequent
# record down by one row (old row 157 becomes row 158).
$ws.Rows.Item(68).Insert()

$ws.Cells.Item(68, 1).Value = 3
$ws.Cells.Item(68, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(68, 3).Value = "Coquimbo"
$ws.Cells.Item(68, 4).Value = 44483
$ws.Cells.Item(68, 5).Value = 5
$ws.Cells.Item(68, 6).Value = 100112001
$ws.Cells.Item(68, 7).Value = "Berenjena"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 50
$ws.Cells.Item(68, 11).Value = 9000
$ws.Cells.Item(68, 12).Value = 9000
$ws.Cells.Item(68, 13).Value = 9000
$ws.Cells.Item(68, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(68, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(68, 16).Value = 150
$ws.Cells.Item(68, 17).Value = 60
$ws.Cells.Item(68, 18).Value = "Hortaliza"
